# Update stock-report figures (qty/value columns, row subtotals and grand
# totals) to match the corrected source data. A handful of product rows
# also had their two stock-batch lines swapped (e.g. rows 192/193), so
# those are expressed here as full value swaps across B/E/F/G.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F8").Value = 65
$ws.Range("G8").Value = 15130.7
$ws.Range("B10").Value = 29894.28
$ws.Range("F22").Value = 60
$ws.Range("G22").Value = 444.6
$ws.Range("F28").Value = 71
$ws.Range("G28").Value = 3166.6
$ws.Range("B32").Value = 13121.82
$ws.Range("F68").Value = 53
$ws.Range("G68").Value = 6101.36
$ws.Range("F70").Value = 28
$ws.Range("G70").Value = 3778.6
$ws.Range("F77").Value = 268
$ws.Range("G77").Value = 12526.32
$ws.Range("F86").Value = 77
$ws.Range("G86").Value = 9661.190000000001
$ws.Range("F89").Value = 2
$ws.Range("G89").Value = 94.59999999999999
$ws.Range("B90").Value = 196197.56
$ws.Range("F115").Value = 225
$ws.Range("G115").Value = 21782.25
$ws.Range("B117").Value = 15834.53
$ws.Range("F144").Value = 1164
$ws.Range("G144").Value = 9835.799999999999
$ws.Range("F145").Value = 584
$ws.Range("G145").Value = 4666.16
$ws.Range("B147").Value = 17280.23
$ws.Range("F153").Value = 28
$ws.Range("G153").Value = 1301.72
$ws.Range("B156").Value = 34512.02
$ws.Range("B192").Value = 64973
$ws.Range("E192").Value = 35.4
$ws.Range("F192").Value = 2
$ws.Range("G192").Value = 66.59999999999999
$ws.Range("B193").Value = 48706
$ws.Range("E193").Value = 39.8
$ws.Range("F193").Value = -144
$ws.Range("G193").Value = -4795.2
$ws.Range("F203").Value = 66
$ws.Range("G203").Value = 1330.56
$ws.Range("F205").Value = 29
$ws.Range("G205").Value = 10937.06
$ws.Range("B216").Value = 46890.33
$ws.Range("F233").Value = 121
$ws.Range("G233").Value = 5764.44
$ws.Range("F255").Value = 596
$ws.Range("G255").Value = 102112.68
$ws.Range("B260").Value = 202610.79
$ws.Range("F280").Value = 143
$ws.Range("G280").Value = 24187.02
$ws.Range("F285").Value = 7
$ws.Range("G285").Value = 195.51
$ws.Range("F288").Value = 46
$ws.Range("G288").Value = 4277.54
$ws.Range("F291").Value = 117
$ws.Range("G291").Value = 5032.17
$ws.Range("F293").Value = 46
$ws.Range("G293").Value = 3234.72
$ws.Range("F294").Value = 42
$ws.Range("G294").Value = 2997.12
$ws.Range("F299").Value = 275
$ws.Range("G299").Value = 39776
$ws.Range("F302").Value = 70
$ws.Range("G302").Value = 14762.3
$ws.Range("F303").Value = 39
$ws.Range("G303").Value = 8224.709999999999
$ws.Range("B304").Value = 189630.57
$ws.Range("F320").Value = 63
$ws.Range("G320").Value = 4324.95
$ws.Range("F321").Value = 46
$ws.Range("G321").Value = 2526.32
$ws.Range("F323").Value = 40
$ws.Range("G323").Value = 4221.6
$ws.Range("F328").Value = 51
$ws.Range("G328").Value = 1897.71
$ws.Range("F329").Value = 29
$ws.Range("G329").Value = 4827.63
$ws.Range("B330").Value = 30360.34
$ws.Range("F339").Value = 5
$ws.Range("G339").Value = 237
$ws.Range("F342").Value = 141
$ws.Range("G342").Value = 4465.47
$ws.Range("F345").Value = 70
$ws.Range("G345").Value = 4298.7
$ws.Range("B346").Value = 27854.16
$ws.Range("B364").Value = 65068
$ws.Range("E364").Value = 13.97
$ws.Range("F364").Value = 63
$ws.Range("G364").Value = 828.45
$ws.Range("B365").Value = 53602
$ws.Range("E365").Value = 15.69
$ws.Range("F365").Value = -231
$ws.Range("G365").Value = -3037.65
$ws.Range("B366").Value = 53263
$ws.Range("E366").Value = 15.29
$ws.Range("F366").Value = -309
$ws.Range("G366").Value = -3958.29
$ws.Range("B367").Value = 65066
$ws.Range("E367").Value = 13.61
$ws.Range("F367").Value = 90
$ws.Range("G367").Value = 1152.9
$ws.Range("B375").Value = 64927
$ws.Range("E375").Value = 17.26
$ws.Range("F375").Value = 106
$ws.Range("G375").Value = 1719.32
$ws.Range("B376").Value = 45718
$ws.Range("E376").Value = 19.38
$ws.Range("F376").Value = -294
$ws.Range("G376").Value = -4768.68
$ws.Range("B382").Value = 64919
$ws.Range("E382").Value = 27.97
$ws.Range("F382").Value = 61
$ws.Range("G382").Value = 1604.3
$ws.Range("B383").Value = 45702
$ws.Range("E383").Value = 31.43
$ws.Range("F383").Value = -215
$ws.Range("G383").Value = -5654.5
$ws.Range("B385").Value = 53595
$ws.Range("E385").Value = 17.61
$ws.Range("F385").Value = -335
$ws.Range("G385").Value = -4934.55
$ws.Range("B386").Value = 65067
$ws.Range("E386").Value = 15.65
$ws.Range("F386").Value = 126
$ws.Range("G386").Value = 1855.98
$ws.Range("F397").Value = 2
$ws.Range("G397").Value = 27.96
$ws.Range("F410").Value = 1
$ws.Range("G410").Value = 32.19
$ws.Range("B411").Value = 7990.15
$ws.Range("F430").Value = 2
$ws.Range("G430").Value = 25.78
$ws.Range("F434").Value = 32
$ws.Range("G434").Value = 1044.48
$ws.Range("B435").Value = 1229.14
$ws.Range("B442").Value = 53319
$ws.Range("E442").Value = 310.64
$ws.Range("F442").Value = -6
$ws.Range("G442").Value = -1643.52
$ws.Range("B443").Value = 64810
$ws.Range("E443").Value = 291.22
$ws.Range("F443").Value = 4
$ws.Range("G443").Value = 1095.68
$ws.Range("B473").Value = 64830
$ws.Range("E473").Value = 34.9
$ws.Range("F473").Value = 108
$ws.Range("G473").Value = 3545.64
$ws.Range("B474").Value = 60022
$ws.Range("E474").Value = 37.22
$ws.Range("F474").Value = -113
$ws.Range("G474").Value = -3709.79
$ws.Range("F539").Value = 46
$ws.Range("G539").Value = 11917.22
$ws.Range("F542").Value = 53
$ws.Range("G542").Value = 6865.09
$ws.Range("B547").Value = 23675.65
$ws.Range("F552").Value = 23
$ws.Range("G552").Value = 2341.17
$ws.Range("F554").Value = 5
$ws.Range("G554").Value = 186.4
$ws.Range("F555").Value = 34
$ws.Range("G555").Value = 2365.04
$ws.Range("B560").Value = 7732.34
$ws.Range("B572").Value = 65079
$ws.Range("F572").Value = 18
$ws.Range("G572").Value = 735.66
$ws.Range("B573").Value = 65362
$ws.Range("F573").Value = 26
$ws.Range("G573").Value = 1062.62
$ws.Range("F578").Value = 92
$ws.Range("G578").Value = 4589.88
$ws.Range("F581").Value = 12
$ws.Range("G581").Value = 2901.6
$ws.Range("F582").Value = 48
$ws.Range("G582").Value = 2735.52
$ws.Range("B583").Value = 25113.22
$ws.Range("F599").Value = 1907
$ws.Range("G599").Value = 311050.77
$ws.Range("F601").Value = 453
$ws.Range("G601").Value = 128140.11
$ws.Range("F602").Value = 345
$ws.Range("G602").Value = 49904.25
$ws.Range("B606").Value = 489943.18
$ws.Range("F612").Value = 34
$ws.Range("G612").Value = 1393.66
$ws.Range("F615").Value = 84
$ws.Range("G615").Value = 10529.4
$ws.Range("B618").Value = 45839.94
$ws.Range("B619").Value = 1914789.1
$ws.Range("B620").Value = 1914789.1
